$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 109.14286
$ws.Range("I8").Value = 109.14286
$ws.Range("K8").Value = 327.42858
$ws.Range("M8").Value = -188.42858
$ws.Range("H17").Value = 5556978.5
$ws.Range("J17").Value = 5883801
$ws.Range("L17").Value = 17651403
$ws.Range("N17").Value = -17651739
$ws.Range("H39").Value = 829.5
$ws.Range("I39").Value = 106
$ws.Range("J39").Value = 3000
$ws.Range("K39").Value = 318
$ws.Range("L39").Value = 9000
$ws.Range("M39").Value = -22
$ws.Range("N39").Value = -9592
$ws.Range("H129").Value = 1253.75
$ws.Range("I129").Value = 1253.75
$ws.Range("K129").Value = 3761.25
$ws.Range("M129").Value = 1238.75
$ws.Range("H132").Value = 5635.48
$ws.Range("I132").Value = 1489.4117
$ws.Range("K132").Value = 4468.2351
$ws.Range("M132").Value = -1938.2351
$ws.Range("H138").Value = 2192.6155
$ws.Range("I138").Value = 1634.75
$ws.Range("J138").Value = 2505.02
$ws.Range("K138").Value = 4904.25
$ws.Range("L138").Value = 7515.059999999999
$ws.Range("M138").Value = 235.75
$ws.Range("N138").Value = -17795.06

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 4307.926
$ws.Range("I110").Value = 2805.682
$ws.Range("K110").Value = 2805.682
$ws.Range("M110").Value = -760.6819999999998
$ws.Range("H122").Value = 3820.739
$ws.Range("I122").Value = 3703.0833
$ws.Range("J122").Value = 3949.0908
$ws.Range("K122").Value = 11109.2499
$ws.Range("L122").Value = 11847.2724
$ws.Range("M122").Value = -8659.249899999999
$ws.Range("N122").Value = -16747.2724
$ws.Range("H132").Value = 3843.4443
$ws.Range("I132").Value = 2765.625
$ws.Range("K132").Value = 8296.875
$ws.Range("M132").Value = -5766.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 762.2308
$ws.Range("I22").Value = 762.2308
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 762.2308
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -589.2308
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H16").Value = 4249.875
$ws.Range("I16").Value = 2263.625
$ws.Range("J16").Value = 6236.125
$ws.Range("K16").Value = 2263.625
$ws.Range("L16").Value = 6236.125
$ws.Range("M16").Value = -1976.625
$ws.Range("N16").Value = -6810.125
$ws.Range("H22").Value = 347.625
$ws.Range("J22").Value = 395.25
$ws.Range("L22").Value = 395.25
$ws.Range("N22").Value = -1095.25
$ws.Range("H33").Value = 566.3333
$ws.Range("I33").Value = 349.5
$ws.Range("J33").Value = 1000
$ws.Range("K33").Value = 349.5
$ws.Range("L33").Value = 1000
$ws.Range("M33").Value = 29.5
$ws.Range("N33").Value = -1758
$ws.Range("H97").Value = 35000
$ws.Range("J97").Value = 35000
$ws.Range("L97").Value = 35000
$ws.Range("N97").Value = -36982
$ws.Range("H113").Value = 4249.875
$ws.Range("I113").Value = 2263.625
$ws.Range("J113").Value = 6236.125
$ws.Range("K113").Value = 2263.625
$ws.Range("L113").Value = 6236.125
$ws.Range("M113").Value = -93.625
$ws.Range("N113").Value = -10576.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8066589.5
$ws.Range("I4").Value = 19187074
$ws.Range("J4").Value = 5054791.5
$ws.Range("K4").Value = 57561222
$ws.Range("L4").Value = 15164374.5
$ws.Range("M4").Value = -57561110
$ws.Range("N4").Value = -15164598.5
$ws.Range("H23").Value = 12722.25
$ws.Range("I23").Value = 276
$ws.Range("K23").Value = 828
$ws.Range("M23").Value = -593
$ws.Range("H28").Value = 3976.3333
$ws.Range("I28").Value = 4114.5
$ws.Range("J28").Value = 3700
$ws.Range("K28").Value = 12343.5
$ws.Range("L28").Value = 11100
$ws.Range("M28").Value = -12111.5
$ws.Range("N28").Value = -11564
$ws.Range("H55").Value = 687.2857
$ws.Range("I55").Value = 634.3333
$ws.Range("J55").Value = 727
$ws.Range("K55").Value = 1902.9999
$ws.Range("L55").Value = 2181
$ws.Range("M55").Value = -1725.9999
$ws.Range("N55").Value = -2535
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H87").Value = 5000
$ws.Range("I87").Value = 5000
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 15000
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -13752
$ws.Range("N87").ClearContents()
$ws.Range("H88").Value = 14999.5
$ws.Range("J88").Value = 14999.5
$ws.Range("L88").Value = 44998.5
$ws.Range("N88").Value = -45854.5
$ws.Range("H90").Value = 5000
$ws.Range("I90").Value = 5000
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 45000
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -38760
$ws.Range("N90").ClearContents()
$ws.Range("H91").Value = 14999.5
$ws.Range("J91").Value = 14999.5
$ws.Range("L91").Value = 44998.5
$ws.Range("N91").Value = -47962.5
$ws.Range("H99").Value = 1150
$ws.Range("I99").Value = 750
$ws.Range("K99").Value = 2250
$ws.Range("M99").Value = -4
$ws.Range("H131").Value = 20835306
$ws.Range("I131").Value = 38462540
$ws.Range("J131").Value = 3123.7273
$ws.Range("K131").Value = 115387620
$ws.Range("L131").Value = 9371.1819
$ws.Range("M131").Value = -115382580
$ws.Range("N131").Value = -19451.1819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3866.6667
$ws.Range("I7").Value = 3866.6667
$ws.Range("K7").Value = 3866.6667
$ws.Range("M7").Value = -3754.6667
$ws.Range("H46").Value = 12194
$ws.Range("I46").Value = 6785.684
$ws.Range("K46").Value = 6785.684
$ws.Range("M46").Value = -6597.684
$ws.Range("H74").Value = 49330.668
$ws.Range("I74").Value = 49330.668
$ws.Range("K74").Value = 49330.668
$ws.Range("M74").Value = -48332.668
$ws.Range("H77").Value = 49330.668
$ws.Range("I77").Value = 49330.668
$ws.Range("K77").Value = 147992.004
$ws.Range("M77").Value = -143000.004
$ws.Range("H122").Value = 4030.7058
$ws.Range("I122").Value = 3954.6
$ws.Range("K122").Value = 11863.8
$ws.Range("M122").Value = -9413.799999999999
$ws.Range("H126").Value = 3866.6667
$ws.Range("I126").Value = 3866.6667
$ws.Range("K126").Value = 11600.0001
$ws.Range("M126").Value = -9130.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 9665
$ws.Range("J15").Value = 9747.5
$ws.Range("L15").Value = 9747.5
$ws.Range("N15").Value = -10323.5
$ws.Range("H100").Value = 414.33334
$ws.Range("I100").Value = 481.16666
$ws.Range("K100").Value = 962.33332
$ws.Range("M100").Value = -421.33332
$ws.Range("H107").Value = 2080
$ws.Range("J107").Value = 2225
$ws.Range("L107").Value = 6675
$ws.Range("N107").Value = -10515
$ws.Range("H122").Value = 4494.759
$ws.Range("J122").Value = 4182.5
$ws.Range("L122").Value = 12547.5
$ws.Range("N122").Value = -17447.5
